$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 62
$ws.Range("H62").Value = 2935
$ws.Range("I62").Value = 2568.3333
$ws.Range("J62").Value = 3118.3333
$ws.Range("K62").Value = 2568.3333
$ws.Range("L62").Value = 3118.3333
$ws.Range("M62").Value = -1944.3333
$ws.Range("N62").Value = -4366.3333

# Row 65
$ws.Range("H65").Value = 2935
$ws.Range("I65").Value = 2568.3333
$ws.Range("J65").Value = 3118.3333
$ws.Range("K65").Value = 12841.6665
$ws.Range("L65").Value = 15591.6665
$ws.Range("M65").Value = -9721.666499999999
$ws.Range("N65").Value = -21831.6665

# Row 98
$ws.Range("H98").Value = 975.3333
$ws.Range("I98").Value = 999.62964
$ws.Range("J98").Value = 866
$ws.Range("K98").Value = 999.62964
$ws.Range("L98").Value = 866
$ws.Range("M98").Value = 498.37036
$ws.Range("N98").Value = -3862

# Row 116
$ws.Range("H116").Value = 27781396
$ws.Range("I116").Value = 3857.8572
$ws.Range("J116").Value = 66669948
$ws.Range("K116").Value = 3857.8572
$ws.Range("L116").Value = 66669948
$ws.Range("M116").Value = -415.8571999999999
$ws.Range("N116").Value = -66676832

# Row 122
$ws.Range("H122").Value = 975.3333
$ws.Range("I122").Value = 999.62964
$ws.Range("J122").Value = 866
$ws.Range("K122").Value = 2998.88892
$ws.Range("L122").Value = 2598
$ws.Range("M122").Value = -548.8889199999999
$ws.Range("N122").Value = -7498

# Row 132
$ws.Range("H132").Value = 1991.9778
$ws.Range("I132").Value = 1720.6897
$ws.Range("J132").Value = 2483.6875
$ws.Range("K132").Value = 5162.0691
$ws.Range("L132").Value = 7451.0625
$ws.Range("M132").Value = -2632.0691
$ws.Range("N132").Value = -12511.0625

# Row 137
$ws.Range("H137").Value = 2380.2952
$ws.Range("I137").Value = 1377.1945
$ws.Range("J137").Value = 3824.76
$ws.Range("K137").Value = 4131.583500000001
$ws.Range("L137").Value = 11474.28
$ws.Range("M137").Value = -1581.583500000001
$ws.Range("N137").Value = -16574.28

# Row 138
$ws.Range("H138").Value = 2725.9485
$ws.Range("I138").Value = 1418.6428
$ws.Range("J138").Value = 3724.2546
$ws.Range("K138").Value = 4255.928400000001
$ws.Range("L138").Value = 11172.7638
$ws.Range("M138").Value = 884.0715999999993
$ws.Range("N138").Value = -21452.7638

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 5959.316
$ws.Range("I32").Value = 5959.316
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5959.316
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -5672.316
$ws.Range("N32").ClearContents()

# Row 54
$ws.Range("H54").Value = 19800
$ws.Range("J54").Value = 19800
$ws.Range("L54").Value = 19800
$ws.Range("N54").Value = -21338

# Row 61
$ws.Range("H61").Value = 4677.0312
$ws.Range("I61").Value = 3452.125
$ws.Range("K61").Value = 3452.125
$ws.Range("M61").Value = -3240.125

# Row 110
$ws.Range("H110").Value = 1860.75
$ws.Range("I110").Value = 918.25
$ws.Range("K110").Value = 918.25
$ws.Range("M110").Value = 1126.75

# Row 132
$ws.Range("H132").Value = 5855.185
$ws.Range("I132").Value = 4151.7095
$ws.Range("J132").Value = 8151.174
$ws.Range("K132").Value = 12455.1285
$ws.Range("L132").Value = 24453.522
$ws.Range("M132").Value = -9925.128499999999
$ws.Range("N132").Value = -29513.522

# Row 136
$ws.Range("H136").Value = 4677.0312
$ws.Range("I136").Value = 3452.125
$ws.Range("K136").Value = 10356.375
$ws.Range("M136").Value = -7806.375

# Row 139
$ws.Range("H139").Value = 72645
$ws.Range("J139").Value = 72645
$ws.Range("L139").Value = 72645
$ws.Range("N139").Value = -82925

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 81
$ws.Range("H81").Value = 16385.555
$ws.Range("J81").Value = 16385.555
$ws.Range("L81").Value = 16385.555
$ws.Range("N81").Value = -18507.555

# Row 84
$ws.Range("H84").Value = 16385.555
$ws.Range("J84").Value = 16385.555
$ws.Range("L84").Value = 49156.665
$ws.Range("N84").Value = -59764.665

# Row 107
$ws.Range("H107").Value = 2657.889
$ws.Range("I107").Value = 2474.4285
$ws.Range("J107").Value = 3300
$ws.Range("K107").Value = 2474.4285
$ws.Range("L107").Value = 3300
$ws.Range("M107").Value = -554.4285
$ws.Range("N107").Value = -7140

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Range("H16").Value = 1161.4286
$ws.Range("J16").Value = 1306
$ws.Range("L16").Value = 1306
$ws.Range("N16").Value = -1880

# Row 58
$ws.Range("H58").Value = 1685926.1
$ws.Range("I58").Value = 2599025
$ws.Range("J58").Value = 3901.7896
$ws.Range("K58").Value = 2599025
$ws.Range("L58").Value = 3901.7896
$ws.Range("M58").Value = -2598822
$ws.Range("N58").Value = -4307.7896

# Row 113
$ws.Range("H113").Value = 1161.4286
$ws.Range("J113").Value = 1306
$ws.Range("L113").Value = 1306
$ws.Range("N113").Value = -5646

# Row 136
$ws.Range("H136").Value = 1685926.1
$ws.Range("I136").Value = 2599025
$ws.Range("J136").Value = 3901.7896
$ws.Range("K136").Value = 7797075
$ws.Range("L136").Value = 11705.3688
$ws.Range("M136").Value = -7794525
$ws.Range("N136").Value = -16805.3688

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 2
$ws.Range("H2").Value = 17.727272
$ws.Range("I2").Value = 26.88889
$ws.Range("J2").Value = 11.384615
$ws.Range("K2").Value = 161.33334
$ws.Range("L2").Value = 68.30769000000001
$ws.Range("M2").Value = -48.33333999999999
$ws.Range("N2").Value = -294.30769

# Row 56
$ws.Range("H56").Value = 71069.87
$ws.Range("I56").Value = 71069.87
$ws.Range("K56").Value = 71069.87
$ws.Range("M56").Value = -70539.87

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 80
$ws.Range("H80").Value = 7885.5
$ws.Range("I80").Value = 11771
$ws.Range("K80").Value = 11771
$ws.Range("M80").Value = -10773

# Row 83
$ws.Range("H83").Value = 7885.5
$ws.Range("I83").Value = 11771
$ws.Range("K83").Value = 58855
$ws.Range("M83").Value = -53863

# Row 132
$ws.Range("H132").Value = 2320.4243
$ws.Range("I132").Value = 1880.091
$ws.Range("J132").Value = 3201.0908
$ws.Range("K132").Value = 5640.272999999999
$ws.Range("L132").Value = 9603.2724
$ws.Range("M132").Value = -3110.272999999999
$ws.Range("N132").Value = -14663.2724

# Row 136
$ws.Range("H136").Value = 6821.35
$ws.Range("J136").Value = 6821.35
$ws.Range("L136").Value = 20464.05
$ws.Range("N136").Value = -25564.05

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 61
$ws.Range("H61").Value = 63778.4
$ws.Range("I61").Value = 100001.336
$ws.Range("J61").Value = 9444
$ws.Range("K61").Value = 100001.336
$ws.Range("L61").Value = 9444
$ws.Range("M61").Value = -99799.336
$ws.Range("N61").Value = -9848

# Row 113
$ws.Range("H113").Value = 63778.4
$ws.Range("I113").Value = 100001.336
$ws.Range("J113").Value = 9444
$ws.Range("K113").Value = 100001.336
$ws.Range("L113").Value = 9444
$ws.Range("M113").Value = -97831.336
$ws.Range("N113").Value = -13784

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 31
$ws.Range("H31").Value = 12000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 12000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 12000
$ws.Range("N31").Value = -12696
$ws.Range("M31").ClearContents()

# Row 113
$ws.Range("H113").Value = 1181.44
$ws.Range("I113").Value = 516.4286
$ws.Range("J113").Value = 2027.8182
$ws.Range("K113").Value = 1549.2858
$ws.Range("L113").Value = 6083.4546
$ws.Range("M113").Value = 620.7142000000001
$ws.Range("N113").Value = -10423.4546

# Row 122
$ws.Range("H122").Value = 3294.9546
$ws.Range("I122").Value = 1832.4445
$ws.Range("K122").Value = 5497.333500000001
$ws.Range("M122").Value = -3047.333500000001

# Row 136
$ws.Range("H136").Value = 3675.5278
$ws.Range("I136").Value = 1571.7174
$ws.Range("J136").Value = 7397.654
$ws.Range("K136").Value = 4715.1522
$ws.Range("L136").Value = 22192.962
$ws.Range("M136").Value = -2165.1522
$ws.Range("N136").Value = -27292.962
